$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1231.1044921875
$ws.Range("C2").Value = 0.9236
$ws.Range("D2").Value = 0.9194999933242798
$ws.Range("E2").Value = 1.359200000762939
$ws.Range("F2").Value = 0.5264999866485596
$ws.Range("H2").Value = 0.793

$ws.Range("B3").Value = 1225.90185546875
$ws.Range("C3").Value = 0.9807
$ws.Range("D3").Value = 0.95
$ws.Range("E3").Value = 1.471500039100647
$ws.Range("F3").Value = 0.6485000252723694
$ws.Range("H3").Value = 1.0635

$ws.Range("B4").Value = 813.3743896484375
$ws.Range("C4").Value = 0.9502
$ws.Range("D4").Value = 0.9083
$ws.Range("E4").Value = 1.560099959373474
$ws.Range("F4").Value = 0.7324000000953674
$ws.Range("H4").Value = 0.6937

$ws.Range("B5").Value = 789.5330200195312
$ws.Range("C5").Value = 0.8328
$ws.Range("D5").Value = 0.8254
$ws.Range("E5").Value = 1.06630003452301
$ws.Range("F5").Value = 0.5223000049591064
$ws.Range("H5").Value = -0.041

$ws.Range("B6").Value = 1101.351684570312
$ws.Range("C6").Value = 0.8713
$ws.Range("D6").Value = 0.867
$ws.Range("E6").Value = 1.078400015830994
$ws.Range("F6").Value = 0.6553999781608582
$ws.Range("H6").Value = 0.3283

$ws.Range("B7").Value = 859.1129760742188
$ws.Range("C7").Value = 0.8652
$ws.Range("D7").Value = 0.8650000095367432
$ws.Range("E7").Value = 1.030500054359436
$ws.Range("F7").Value = 0.723800003528595
$ws.Range("H7").Value = 0.3104

$ws.Range("B8").Value = 947.958984375
$ws.Range("C8").Value = 0.8494
$ws.Range("D8").Value = 0.8488
$ws.Range("E8").Value = 1.085000038146973
$ws.Range("F8").Value = 0.7342000007629395
$ws.Range("H8").Value = 0.1665

$ws.Range("B9").Value = 6968.33740234375
$ws.Range("C9").Value = 0.898
$ws.Range("D9").Value = 0.8848
$ws.Range("E9").Value = 1.560099959373474
$ws.Range("F9").Value = 0.5223000049591064
$ws.Range("H9").Value = 3.3144
